$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E3").Value = 10.33

$ws.Range("E4").Value = 10.03
$ws.Range("F4").Value = 9.98
$ws.Range("G4").Value = 10.29

$ws.Range("C5").Value = 9.67
$ws.Range("D5").Value = 9.970000000000001
$ws.Range("F5").Value = 10.26
$ws.Range("G5").Value = 9.130000000000001
$ws.Range("J5").Value = 6.4

$ws.Range("D6").Value = 10.02
$ws.Range("E6").Value = 9.74
$ws.Range("G6").Value = 10.44
$ws.Range("I6").Value = 10.18

$ws.Range("D7").Value = 9.710000000000001
$ws.Range("E7").Value = 10.87
$ws.Range("F7").Value = 9.56

$ws.Range("J8").Value = 12.4

$ws.Range("F9").Value = 9.82

$ws.Range("E10").Value = 13.6
$ws.Range("H10").Value = 7.6
